$d = $word.ActiveDocument

$styleNames = @(
    "Heading 1 Char",
    "Heading 2 Char",
    "Heading 3 Char",
    "Heading 4 Char",
    "Heading 5 Char",
    "Heading 6 Char",
    "Title Char",
    "Comment Text Char"
)

foreach ($name in $styleNames) {
    $s = $d.Styles($name)
    $f = $s.Font
    if ($f.Name -eq "Bookerly") {
        $f.Name = "Times New Roman"
    }
    if ($f.NameAscii -eq "Bookerly") {
        $f.NameAscii = "Times New Roman"
    }
}
